$d = $word.ActiveDocument

$replacements = @(
    @("152÷2=", "159÷7="),
    @("303÷2=", "109÷4="),
    @("485÷8=", "361÷2="),
    @("752÷8=", "534÷9="),
    @("429÷9=", "371÷5="),
    @("170÷7=", "853÷2="),
    @("388÷7=", "265÷4="),
    @("291÷2=", "246÷6="),
    @("869÷4=", "134÷5="),
    @("494÷2=", "912÷9="),
    @("586÷6=", "947÷3="),
    @("776÷6=", "331÷4="),
    @("612÷3=", "414÷9="),
    @("711÷9=", "364÷7="),
    @("318÷4=", "940÷2="),
    @("155÷6=", "816÷5="),
    @("135÷8=", "534÷2="),
    @("848÷8=", "873÷3="),
    @("762÷8=", "899÷8="),
    @("145÷2=", "400÷7="),
    @("281÷7=", "271÷4="),
    @("757÷5=", "547÷9="),
    @("512÷5=", "581÷8="),
    @("995÷5=", "287÷5="),
    @("407÷9=", "613÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($replacements.Count) division expressions."
